# Applies the "atualização atividade aula 03" edit described by the diff.
#
# Notes on the interop engine's behaviour (learned by experimentation):
#  - A paragraph object fetched via $d.Paragraphs.Item($i) effectively
#    re-resolves "the paragraph currently at position $i" whenever its
#    properties are touched, rather than keeping a fixed reference to the
#    paragraph it originally pointed to. So after inserting/deleting
#    paragraphs earlier in the document, previously-fetched objects for
#    later positions can start referring to the wrong paragraph.
#    To stay safe we always process the document from the BOTTOM up, so
#    edits never shift the index of a paragraph we still need to reach.
#  - Setting $paragraph.Range.Text = "..." only overwrites up to the end
#    of the *first run* in that paragraph, leaving any further runs
#    untouched, when the paragraph has several runs. To really replace the
#    whole paragraph (and collapse it to one run) the range has to be
#    built explicitly from Start to End-1 (excluding the paragraph mark).
#  - For paragraphs where only a sub-string changes while the rest of the
#    text (and its run split) should stay as-is, Find/Execute replace is
#    used instead, since it only touches the runs it actually spans.

function Set-ParagraphText($doc, $index, $text) {
    $p = $doc.Paragraphs.Item($index)
    $startPos = $p.Range.Start
    $endPos = $p.Range.End
    $rng = $doc.Range($startPos, $endPos - 1)
    $rng.Text = $text
}

function Insert-ParagraphsAfter($doc, $index, $texts) {
    # Inserts a sequence of new paragraphs, in order, right after the
    # paragraph currently at $index. Empty strings are left as bare empty
    # paragraphs (no point writing an explicit empty <w:t/> run).
    $cur = $index
    foreach ($t in $texts) {
        $doc.Paragraphs.Item($cur).Range.InsertParagraphAfter()
        $cur = $cur + 1
        if ($t -ne "") {
            $doc.Paragraphs.Item($cur).Range.Text = $t
        }
    }
}

$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Biblioteca block (originally paragraphs 22-28)
# ----------------------------------------------------------------------

# 27: "Biblioteca está organizada em setores" -> removed
$d.Paragraphs.Item(27).Range.Delete()

# 26: "Biblioteca possui usuário" -> "Usuário reserva livro" + new empty para
Set-ParagraphText $d 26 "Usuário reserva livro"
$d.Paragraphs.Item(26).Range.InsertParagraphAfter()

# 24: "Biblioteca tem funcionário" -> "Funcionário (...)" + new paragraphs
Set-ParagraphText $d 24 "Funcionário (matrícula, nome, CPF, endereço, telefone, e-mail)"
Insert-ParagraphsAfter $d 24 @(
    "Usuário (código, nome, CPF, endereço, telefone, e-mail)",
    "Livro (ISBN, título, autor, ano de publicação, edição, volume)",
    "Biblioteca possui funcionário",
    "Biblioteca tem usuário"
)

# 23: "Biblioteca (...)" definition line, drop the extra attribute list
# (only the inner attribute-list text changes; leading/trailing runs stay)
$d.Content.Find.Execute(", funcionário, usuário, livro, setores", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ----------------------------------------------------------------------
# Escola block (originally paragraphs 15-21)
# ----------------------------------------------------------------------

# 20: "Escola possui cursos" -> removed
$d.Paragraphs.Item(20).Range.Delete()
# 19: "Escola tem aluno" -> removed
$d.Paragraphs.Item(19).Range.Delete()
# 18: "Escola tem professor" -> removed
$d.Paragraphs.Item(18).Range.Delete()

# 17: "Escola possui diretor" -> "Professor (...)" + new paragraphs
Set-ParagraphText $d 17 "Professor (matrícula, nome, telefone, e-mail, endereço)"
Insert-ParagraphsAfter $d 17 @(
    "Aluno (matrícula, nome, telefone, e-mail, endereço)",
    "Disciplina (código, nome, área do conhecimento, duração)",
    "Escola possui professor",
    "Escola possui aluno",
    "Professor ensina aluno",
    "Professor leciona disciplina",
    "Aluno cursa disciplina",
    ""
)

# 16: "Escola (...)" definition line, drop the extra attribute list
$d.Content.Find.Execute("escola, nome, CNPJ, endereço, telefone, diretor, professor, aluno, coordenador, cursos", $true, $false, $false, $false, $false, $true, 1, $false, "escola, nome, CNPJ, endereço, telefone", 2) | Out-Null

# ----------------------------------------------------------------------
# Loja de roupas block (originally paragraphs 8-14)
# ----------------------------------------------------------------------

# 13: "Loja aceita forma de pagamento" -> "Loja de roupas tem cliente" + new paragraphs
Set-ParagraphText $d 13 "Loja de roupas tem cliente"
Insert-ParagraphsAfter $d 13 @(
    "Loja de roupas compra de fornecedor",
    "Loja de roupa vende produto",
    "Fornecedor tem vendedor",
    "Vendedor vende produto",
    "Cliente compra produto"
)

# 11: "Loja de roupas possui gerente" -> "Funcionário (...)" + new paragraphs
Set-ParagraphText $d 11 "Funcionário (matrícula, nome, endereço, telefone, e-mail)"
Insert-ParagraphsAfter $d 11 @(
    "Cliente (código, nome, CPF/CNPJ, endereço, telefone, e-mail)",
    "Fornecedor (código, nome, CNPJ, endereço, telefone, e-mail)",
    "Vendedor (matrícula, nome, CPF, endereço, telefone, e-mail)",
    "Produto (código, nome, marca, quantidade em estoque, valor unitário, valor total)"
)

# 10: "Loja de roupas (...)" definition line, drop the extra attribute list
$d.Content.Find.Execute("gerente, telefone, funcionário, forma de pagamento, departamento", $true, $false, $false, $false, $false, $true, 1, $false, "telefone", 2) | Out-Null

# 8: "Crie o modelo conceitual dos seguintes cenários:" -> add blank paragraph after it
$d.Paragraphs.Item(8).Range.InsertParagraphAfter()

Write-Host "done"
